$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove (clear) the two now-unneeded commission test-case rows (111 and 118).
# Clearing the full data span (C:I) empties every cell in the row so the row
# itself disappears from the saved sheet, while all the other row numbers
# stay exactly where they are (no shifting of the remaining data).
$ws.Range("C111:I111").ClearContents() | Out-Null
$ws.Range("C118:I118").ClearContents() | Out-Null

# Move / record the active selection as it was left after the edit.
$ws.Range("G19").Select() | Out-Null
